# Updates the "Metadata" sheet of the Engagement Product ValueSet workbook:
#  - Version bump 5.0.0 -> 6.0.0
#  - Date bump to the new publication timestamp
#  - Publisher's Contact/ContactDetail rows replaced with Publisher = "Alvearie Team"
#    and a new "Jurisdiction" = "United States of America" row
#  - Removes the now-redundant duplicate row, shrinking the sheet from 15 to 14 rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Drop the trailing duplicate row so the table is back to 14 data rows.
$ws.Rows.Item(15).Delete()

# Rewrite the Property/Value pairs top-to-bottom with the new content.
$ws.Cells.Item(1, 1).Value = "Property"
$ws.Cells.Item(1, 2).Value = "Value"

$ws.Cells.Item(2, 1).Value = "URL"
$ws.Cells.Item(2, 2).Value = "http://ibm.com/fhir/cdm/ValueSet/eng-product"

$ws.Cells.Item(3, 1).Value = "Version"
$ws.Cells.Item(3, 2).Value = "6.0.0"

$ws.Cells.Item(4, 1).Value = "Name"
$ws.Cells.Item(4, 2).Value = "EngagementProductValueSet"

$ws.Cells.Item(5, 1).Value = "Title"
$ws.Cells.Item(5, 2).Value = "Engagement Product Value Set"

$ws.Cells.Item(6, 1).Value = "Status"
$ws.Cells.Item(6, 2).Value = "active"

$ws.Cells.Item(7, 1).Value = "Experimental"
$ws.Cells.Item(7, 2).Value = ""

$ws.Cells.Item(8, 1).Value = "Date"
$ws.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

$ws.Cells.Item(9, 1).Value = "Publisher"
$ws.Cells.Item(9, 2).Value = "Alvearie Team"

$ws.Cells.Item(10, 1).Value = "Jurisdiction"
$ws.Cells.Item(10, 2).Value = "United States of America"

$ws.Cells.Item(11, 1).Value = "Description"
$ws.Cells.Item(11, 2).Value = "IBM Watson Health engagement products"

$ws.Cells.Item(12, 1).Value = "Purpose"
$ws.Cells.Item(12, 2).Value = ""

$ws.Cells.Item(13, 1).Value = "Copyright"
$ws.Cells.Item(13, 2).Value = ""

$ws.Cells.Item(14, 1).Value = "Immutable"
$ws.Cells.Item(14, 2).Value = "BooleanType[null]"
